# Scheduled market-price refresh: update computed Leve profit columns (H:N)
# across the per-job worksheets, matching the latest currentAveragePrice pulls.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 138.33333
$ws.Cells.Item(33, 9).Value = 134.28572
$ws.Cells.Item(33, 11).Value = 134.28572
$ws.Cells.Item(33, 13).Value = 94.71428
$ws.Cells.Item(38, 8).Value = 1649.0834
$ws.Cells.Item(38, 9).Value = 163.33333
$ws.Cells.Item(38, 10).Value = 3134.8333
$ws.Cells.Item(38, 11).Value = 489.99999
$ws.Cells.Item(38, 12).Value = 9404.499899999999
$ws.Cells.Item(38, 13).Value = -117.99999
$ws.Cells.Item(38, 14).Value = -10148.4999
$ws.Cells.Item(76, 8).Value = 7410600
$ws.Cells.Item(76, 9).Value = 10103910
$ws.Cells.Item(76, 11).Value = 10103910
$ws.Cells.Item(76, 13).Value = -10103595
$ws.Cells.Item(79, 8).Value = 7410600
$ws.Cells.Item(79, 9).Value = 10103910
$ws.Cells.Item(79, 11).Value = 10103910
$ws.Cells.Item(79, 13).Value = -10102818
$ws.Cells.Item(137, 8).Value = 19231742
$ws.Cells.Item(137, 9).Value = 25641726
$ws.Cells.Item(137, 10).Value = 1788.8462
$ws.Cells.Item(137, 11).Value = 76925178
$ws.Cells.Item(137, 12).Value = 5366.5386
$ws.Cells.Item(137, 13).Value = -76922628
$ws.Cells.Item(137, 14).Value = -10466.5386
$ws.Cells.Item(141, 8).Value = 1958.7097
$ws.Cells.Item(141, 9).Value = 1191.3334
$ws.Cells.Item(141, 10).Value = 5156.1113
$ws.Cells.Item(141, 11).Value = 3574.0002
$ws.Cells.Item(141, 12).Value = 15468.3339
$ws.Cells.Item(141, 13).Value = 1605.9998
$ws.Cells.Item(141, 14).Value = -25828.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 5633.3335
$ws.Cells.Item(31, 9).Value = 5633.3335
$ws.Cells.Item(31, 11).Value = 5633.3335
$ws.Cells.Item(31, 13).Value = -5339.3335
$ws.Cells.Item(32, 8).Value = 22669.81
$ws.Cells.Item(32, 9).Value = 3662.4082
$ws.Cells.Item(32, 11).Value = 3662.4082
$ws.Cells.Item(32, 13).Value = -3375.4082
$ws.Cells.Item(35, 8).Value = 1950
$ws.Cells.Item(35, 9).Value = 1950
$ws.Cells.Item(35, 11).Value = 1950
$ws.Cells.Item(35, 13).Value = -1544
$ws.Cells.Item(88, 8).Value = 5479.7
$ws.Cells.Item(88, 9).Value = 3180
$ws.Cells.Item(88, 10).Value = 7779.4
$ws.Cells.Item(88, 11).Value = 3180
$ws.Cells.Item(88, 12).Value = 7779.4
$ws.Cells.Item(88, 13).Value = -2774
$ws.Cells.Item(88, 14).Value = -8591.4
$ws.Cells.Item(91, 8).Value = 5479.7
$ws.Cells.Item(91, 9).Value = 3180
$ws.Cells.Item(91, 10).Value = 7779.4
$ws.Cells.Item(91, 11).Value = 3180
$ws.Cells.Item(91, 12).Value = 7779.4
$ws.Cells.Item(91, 13).Value = -1776
$ws.Cells.Item(91, 14).Value = -10587.4
$ws.Cells.Item(122, 8).Value = 5947.826
$ws.Cells.Item(122, 9).Value = 6576.15
$ws.Cells.Item(122, 11).Value = 19728.45
$ws.Cells.Item(122, 13).Value = -17278.45
$ws.Cells.Item(132, 8).Value = 2081.8906
$ws.Cells.Item(132, 9).Value = 1812.8628
$ws.Cells.Item(132, 10).Value = 3137.3076
$ws.Cells.Item(132, 11).Value = 5438.588400000001
$ws.Cells.Item(132, 12).Value = 9411.9228
$ws.Cells.Item(132, 13).Value = -2908.588400000001
$ws.Cells.Item(132, 14).Value = -14471.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(37, 8).Value = 11204.167
$ws.Cells.Item(37, 9).Value = 3445
$ws.Cells.Item(37, 11).Value = 3445
$ws.Cells.Item(37, 13).Value = -3308
$ws.Cells.Item(86, 8).Value = 10665.333
$ws.Cells.Item(86, 9).Value = 4666.6665
$ws.Cells.Item(86, 10).Value = 12664.889
$ws.Cells.Item(86, 11).Value = 4666.6665
$ws.Cells.Item(86, 12).Value = 12664.889
$ws.Cells.Item(86, 13).Value = -3543.6665
$ws.Cells.Item(86, 14).Value = -14910.889
$ws.Cells.Item(89, 8).Value = 10665.333
$ws.Cells.Item(89, 9).Value = 4666.6665
$ws.Cells.Item(89, 10).Value = 12664.889
$ws.Cells.Item(89, 11).Value = 23333.3325
$ws.Cells.Item(89, 12).Value = 63324.44499999999
$ws.Cells.Item(89, 13).Value = -17717.3325
$ws.Cells.Item(89, 14).Value = -74556.44499999999
$ws.Cells.Item(94, 8).Value = 1035.9429
$ws.Cells.Item(94, 9).Value = 923.1429000000001
$ws.Cells.Item(94, 10).Value = 1487.1428
$ws.Cells.Item(94, 11).Value = 923.1429000000001
$ws.Cells.Item(94, 12).Value = 1487.1428
$ws.Cells.Item(94, 13).Value = -472.1429000000001
$ws.Cells.Item(94, 14).Value = -2389.1428
$ws.Cells.Item(102, 8).Value = 23278
$ws.Cells.Item(102, 9).Value = 12556
$ws.Cells.Item(102, 10).Value = 34000
$ws.Cells.Item(102, 11).Value = 12556
$ws.Cells.Item(102, 12).Value = 34000
$ws.Cells.Item(102, 13).Value = -9311
$ws.Cells.Item(102, 14).Value = -40490

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 60.666668
$ws.Cells.Item(7, 9).Value = 60
$ws.Cells.Item(7, 10).Value = 61.333332
$ws.Cells.Item(7, 11).Value = 60
$ws.Cells.Item(7, 12).Value = 61.333332
$ws.Cells.Item(7, 13).Value = 53
$ws.Cells.Item(7, 14).Value = -287.333332
$ws.Cells.Item(17, 8).Value = 27500
$ws.Cells.Item(17, 9).Value = 50000
$ws.Cells.Item(17, 10).Value = 5000
$ws.Cells.Item(17, 11).Value = 50000
$ws.Cells.Item(17, 12).Value = 5000
$ws.Cells.Item(17, 13).Value = -49826
$ws.Cells.Item(17, 14).Value = -5348
$ws.Cells.Item(20, 8).Value = 49800
$ws.Cells.Item(20, 10).Value = 49800
$ws.Cells.Item(20, 12).Value = 49800
$ws.Cells.Item(20, 14).Value = -50272
$ws.Cells.Item(25, 8).Value = 6500
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 6500
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 6500
$ws.Cells.Item(25, 13).ClearContents()
$ws.Cells.Item(25, 14).Value = -6848
$ws.Cells.Item(30, 8).Value = 49800
$ws.Cells.Item(30, 10).Value = 49800
$ws.Cells.Item(30, 12).Value = 49800
$ws.Cells.Item(30, 14).Value = -49982
$ws.Cells.Item(31, 8).Value = 1892.82
$ws.Cells.Item(31, 9).Value = 1156.4166
$ws.Cells.Item(31, 10).Value = 3786.4285
$ws.Cells.Item(31, 11).Value = 1156.4166
$ws.Cells.Item(31, 12).Value = 3786.4285
$ws.Cells.Item(31, 13).Value = -861.4166
$ws.Cells.Item(31, 14).Value = -4376.4285
$ws.Cells.Item(34, 8).Value = 1892.82
$ws.Cells.Item(34, 9).Value = 1156.4166
$ws.Cells.Item(34, 10).Value = 3786.4285
$ws.Cells.Item(34, 11).Value = 1156.4166
$ws.Cells.Item(34, 12).Value = 3786.4285
$ws.Cells.Item(34, 13).Value = -954.4166
$ws.Cells.Item(34, 14).Value = -4190.4285
$ws.Cells.Item(39, 8).Value = 550
$ws.Cells.Item(39, 9).Value = 550
$ws.Cells.Item(39, 11).Value = 550
$ws.Cells.Item(39, 13).Value = -159
$ws.Cells.Item(41, 8).Value = 12870.75
$ws.Cells.Item(41, 9).Value = 7472.6665
$ws.Cells.Item(41, 10).Value = 29065
$ws.Cells.Item(41, 11).Value = 7472.6665
$ws.Cells.Item(41, 12).Value = 29065
$ws.Cells.Item(41, 14).Value = -29921
$ws.Cells.Item(41, 13).Value = -7044.6665
$ws.Cells.Item(49, 8).Value = 550
$ws.Cells.Item(49, 9).Value = 550
$ws.Cells.Item(49, 11).Value = 550
$ws.Cells.Item(49, 13).Value = -368
$ws.Cells.Item(52, 8).Value = 40000
$ws.Cells.Item(52, 10).Value = 40000
$ws.Cells.Item(52, 12).Value = 40000
$ws.Cells.Item(52, 14).Value = -40588
$ws.Cells.Item(55, 8).Value = 5054.6
$ws.Cells.Item(55, 9).Value = 1536.5
$ws.Cells.Item(55, 11).Value = 1536.5
$ws.Cells.Item(55, 13).Value = -1221.5
$ws.Cells.Item(58, 8).Value = 1469.5333
$ws.Cells.Item(58, 9).Value = 929.62164
$ws.Cells.Item(58, 11).Value = 929.62164
$ws.Cells.Item(58, 13).Value = -726.62164
$ws.Cells.Item(128, 8).Value = 49800
$ws.Cells.Item(128, 10).Value = 49800
$ws.Cells.Item(128, 12).Value = 49800
$ws.Cells.Item(128, 14).Value = -59760
$ws.Cells.Item(130, 8).Value = 48500
$ws.Cells.Item(130, 10).Value = 48500
$ws.Cells.Item(130, 12).Value = 48500
$ws.Cells.Item(130, 14).Value = -58540
$ws.Cells.Item(132, 8).Value = 1876.05
$ws.Cells.Item(132, 9).Value = 1627.2162
$ws.Cells.Item(132, 10).Value = 4945
$ws.Cells.Item(132, 11).Value = 4881.6486
$ws.Cells.Item(132, 12).Value = 14835
$ws.Cells.Item(132, 13).Value = -2351.6486
$ws.Cells.Item(132, 14).Value = -19895
$ws.Cells.Item(136, 8).Value = 1469.5333
$ws.Cells.Item(136, 9).Value = 929.62164
$ws.Cells.Item(136, 11).Value = 2788.86492
$ws.Cells.Item(136, 13).Value = -238.86492

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 2350.258
$ws.Cells.Item(131, 10).Value = 2411.9333
$ws.Cells.Item(131, 12).Value = 7235.7999
$ws.Cells.Item(131, 14).Value = -17315.7999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2055
$ws.Cells.Item(126, 9).Value = 1566.1538
$ws.Cells.Item(126, 10).Value = 2389.4736
$ws.Cells.Item(126, 11).Value = 4698.4614
$ws.Cells.Item(126, 12).Value = 7168.4208
$ws.Cells.Item(126, 13).Value = -2228.4614
$ws.Cells.Item(126, 14).Value = -12108.4208

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1533.3334
$ws.Cells.Item(46, 10).Value = 1800
$ws.Cells.Item(46, 12).Value = 1800
$ws.Cells.Item(46, 14).Value = -2176
$ws.Cells.Item(136, 8).Value = 3533.9805
$ws.Cells.Item(136, 9).Value = 1917.0834
$ws.Cells.Item(136, 11).Value = 5751.2502
$ws.Cells.Item(136, 13).Value = -3201.2502

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 10011500
$ws.Cells.Item(14, 9).Value = 25750
$ws.Cells.Item(14, 11).Value = 25750
$ws.Cells.Item(14, 13).Value = -25582
$ws.Cells.Item(17, 8).Value = 3849.4546
$ws.Cells.Item(17, 9).Value = 4355.5
$ws.Cells.Item(17, 10).Value = 2500
$ws.Cells.Item(17, 11).Value = 4355.5
$ws.Cells.Item(17, 12).Value = 2500
$ws.Cells.Item(17, 13).Value = -4183.5
$ws.Cells.Item(17, 14).Value = -2844
$ws.Cells.Item(132, 8).Value = 8930758
$ws.Cells.Item(132, 9).Value = 14287917
$ws.Cells.Item(132, 10).Value = 2158.2856
$ws.Cells.Item(132, 11).Value = 42863751
$ws.Cells.Item(132, 12).Value = 6474.8568
$ws.Cells.Item(132, 13).Value = -42861221
$ws.Cells.Item(132, 14).Value = -11534.8568
